$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the "time_taken" column (F) on the "data" sheet with refreshed timestamps ---
$newTimestamps = @(
    "2021-10-05 14:19:46.548505",
    "2021-10-05 14:19:46.548514",
    "2021-10-05 14:19:46.548517",
    "2021-10-05 14:19:46.548520",
    "2021-10-05 14:19:46.548523",
    "2021-10-05 14:19:46.548526",
    "2021-10-05 14:19:46.548529",
    "2021-10-05 14:19:46.548531",
    "2021-10-05 14:19:46.548534",
    "2021-10-05 14:19:46.548537",
    "2021-10-05 14:19:46.548540",
    "2021-10-05 14:19:46.548543",
    "2021-10-05 14:19:46.548545",
    "2021-10-05 14:19:46.548548",
    "2021-10-05 14:19:46.548551",
    "2021-10-05 14:19:46.548554",
    "2021-10-05 14:19:46.548557",
    "2021-10-05 14:19:46.548560",
    "2021-10-05 14:19:46.548563",
    "2021-10-05 14:19:46.548565",
    "2021-10-05 14:19:46.548568",
    "2021-10-05 14:19:46.548571",
    "2021-10-05 14:19:46.548574",
    "2021-10-05 14:19:46.548576",
    "2021-10-05 14:19:46.548579",
    "2021-10-05 14:19:46.548582",
    "2021-10-05 14:19:46.548585",
    "2021-10-05 14:19:46.548588",
    "2021-10-05 14:19:46.548590",
    "2021-10-05 14:19:46.548593",
    "2021-10-05 14:19:46.548596",
    "2021-10-05 14:19:46.548598",
    "2021-10-05 14:19:46.548601",
    "2021-10-05 14:19:46.548604",
    "2021-10-05 14:19:46.548607",
    "2021-10-05 14:19:46.548610",
    "2021-10-05 14:19:46.548612",
    "2021-10-05 14:19:46.548615",
    "2021-10-05 14:19:46.548618",
    "2021-10-05 14:19:46.548620",
    "2021-10-05 14:19:46.548623",
    "2021-10-05 14:19:46.548626"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- Add the new "metadata" sheet right after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Copy the bold/bordered header style from the "data" sheet so no new style entries are created
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)

# Data row
$metaSheet.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Corneal abnormalities"
$metaSheet.Range("C2").Value = 250

# "1.9" must stay text (not be coerced to the number 1.9): stamp it as Text,
# write it, then drop the formatting again so no extra style is left behind.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.9"
$metaSheet.Range("D2").ClearFormats()

$metaSheet.Range("E2").Value = "2021-08-31T15:03:20.149221Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:19:46.545057"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/250/?format=json"

$excel.CutCopyMode = 0
$dataSheet.Activate()
